$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells remain text so numeric-looking strings (leading/trailing zeros,
# percent signs) are preserved exactly as in the source diff.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "307.14"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.18%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "35.88"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.46%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.111"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.00%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08083"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.96%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.948"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.07%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.197"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "3.64%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "7.746"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.08%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9276"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.71%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1378"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "12.43%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1896"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2.18%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09189"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-2.90%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03408"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-4.86%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09826"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.16%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001430"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "3.09%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005756"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.88%"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "3.67%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3449"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.19%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "4.59%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.902"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-3.03%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2446"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.76%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04425"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-2.15%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001223"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.54%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004827"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.37%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001243"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.54%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02021"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "4.48%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04916"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "3.44%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007710"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.49%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01011"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "5.96%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1377"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "3.47%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002105"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-0.22%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01158"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "3.99%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006453"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "2.56%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.26%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "64.67"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "0.29%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-19.80%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002105"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.26%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002005"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.26%"
